$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 83, shifting the existing data (rows 83:124)
# down to rows 84:125.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with a new weekly price record.
$ws.Range("A83").Value = 10
$ws.Range("B83").Value = "Vega Modelo de Temuco"
$ws.Range("C83").Value = "La Araucanía"
$ws.Range("D83").Value = 44455
$ws.Range("E83").Value = 9
$ws.Range("F83").Value = 100112005
$ws.Range("G83").Value = "Puerro"
$ws.Range("H83").Value = "Azul de Maquehue"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 40
$ws.Range("K83").Value = 8000
$ws.Range("L83").Value = 8000
$ws.Range("M83").Value = 8000
$ws.Range("N83").Value = "$/docena de paquetes"
$ws.Range("O83").Value = "Provincia de Cautín"
$ws.Range("P83").Value = 667
$ws.Range("Q83").Value = 12
$ws.Range("R83").Value = "Hortaliza"
